$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-parsed as numbers by Excel, to preserve them as text (matching source).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated cell values (Coin, Link, Price, Volume columns).
$ws.Range("D2").Value = "35.141.60"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.814.21"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  +0.73%  "
$ws.Range("D5").Value = "233.35"
$ws.Range("E5").Value = "  +2.19%  "
$ws.Range("D6").Value = "0.613"
$ws.Range("E6").Value = "  +0.21%  "
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").Value = "40.68"
$ws.Range("E8").Value = "  -5.63%  "
$ws.Range("D9").Value = "0.324"
$ws.Range("E9").Value = "  +6.05%  "
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").Value = "0.0997"
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "2.077.74"
$ws.Range("E12").Value = "  -1.59%  "
$ws.Range("D13").Value = "1.829.78"
$ws.Range("E13").Value = "  -0.39%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "0.663"
$ws.Range("E14").Value = "  +0.20%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "11.06"
$ws.Range("E15").Value = "  -4.94%  "
$ws.Range("D16").Value = "4.66"
$ws.Range("E16").Value = "  -1.72%  "
$ws.Range("D17").Value = "35.092.62"
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").Value = "69.60"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "238.70"
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("D21").Value = "11.90"
$ws.Range("E21").Value = "  -1.89%  "
$ws.Range("D22").Value = "4.69"
$ws.Range("E22").Value = "  -0.73%  "
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("E24").Value = "  +3.20%  "
$ws.Range("D25").Value = "172.68"
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "7.83"
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("D27").Value = "17.49"
$ws.Range("E27").Value = "  -2.44%  "
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("D29").Value = "1.61"
$ws.Range("E29").Value = "  +21.28%  "
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "4.18"
$ws.Range("E31").Value = "  +5.82%  "
$ws.Range("D32").Value = "3.329.75"
$ws.Range("E32").Value = "  -5.56%  "
$ws.Range("E33").Value = "  +2.59%  "
$ws.Range("D34").Value = "4.01"
$ws.Range("E34").Value = "  -1.46%  "
$ws.Range("E35").Value = "  -6.20%  "
$ws.Range("E36").Value = "  +5.26%  "
$ws.Range("D37").Value = "92.46"
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("E39").Value = "  -0.20%  "
$ws.Range("B40").Value = "WEMIXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").Value = "1.28"
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "1.311.09"
$ws.Range("E41").Value = "  -2.47%  "
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").Value = "14.53"
$ws.Range("E44").Value = "  -2.16%  "
$ws.Range("D45").Value = "2.29"
$ws.Range("E45").Value = "  -6.00%  "
$ws.Range("E46").Value = "  -2.45%  "
$ws.Range("D47").Value = "6.32"
$ws.Range("E47").Value = "  +4.22%  "
$ws.Range("D48").Value = "0.0511"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("D49").Value = "1.991.98"
$ws.Range("E49").Value = "  -0.95%  "
$ws.Range("E50").Value = "  +0.70%  "
$ws.Range("D51").Value = "0.0653"
$ws.Range("E51").Value = "  +5.13%  "
